$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update the "last updated" timestamp string (row 1) ---
$ws.Range("A1").Value = "Datos actualizados a 6 de Septiembre de 2020 a las 13:19"

# --- Country re-ranking: a few rows swapped places because their totals
#     crossed over after the data refresh. Swap the country-name labels in
#     column A for those row pairs; the numeric columns for every affected
#     row (old and new) are corrected below via the main data table. ---
$ws.Range("A44").Value = "Emiratos Arabes Unidos"
$ws.Range("A45").Value = "Paises Bajos"

$ws.Range("A82").Value = "Libia"
$ws.Range("A83").Value = "Dinamarca"

$ws.Range("A128").Value = "Eslovenia"
$ws.Range("A129").Value = "Gambia"

# --- Refreshed case-count data for the affected rows ---
$values = @{
    "B4" = 6432103
    "C4" = 951
    "D4" = 3707128
    "E4" = 2532125
    "G4" = 32
    "H4" = 192850
    "B6" = 4118583
    "C6" = 7744
    "E6" = 866872
    "G6" = 33
    "H6" = 70712
    "B15" = 386658
    "C15" = 1992
    "D15" = 333900
    "E15" = 30465
    "G15" = 139
    "H15" = 22293
    "B17" = 325157
    "C17" = 1592
    "D17" = 221275
    "E17" = 99403
    "G17" = 32
    "H17" = 4479
    "B26" = 194109
    "C26" = 3444
    "D26" = 138575
    "E26" = 47509
    "G26" = 85
    "H26" = 8025
    "B31" = 120095
    "C31" = 231
    "D31" = 116998
    "E31" = 2894
    "G31" = 1
    "H31" = 203
    "B37" = 95014
    "C37" = 1150
    "D37" = 40307
    "E37" = 50814
    "G37" = 43
    "H37" = 3893
    "B38" = 89582
    "C38" = 619
    "D38" = 80521
    "E38" = 8517
    "G38" = 4
    "H38" = 544
    "B44" = 73984
    "C44" = 513
    "D44" = 66095
    "E44" = 7501
    "H44" = 388
    "B45" = 73862
    "D45" = 0
    "E45" = 0
    "H45" = 6241
    "E55" = 3928
    "G55" = 2
    "H55" = 198
    "B63" = 44401
    "C63" = 444
    "E63" = 5288
    "B71" = 29271
    "C71" = 184
    "D71" = 25043
    "E71" = 3492
    "G71" = 1
    "H71" = 736
    "D74" = 22462
    "E74" = 3063
    "B82" = 17749
    "C82" = 655
    "D82" = 2081
    "E82" = 15383
    "G82" = 13
    "H82" = 285
    "B83" = 17736
    "D83" = 15671
    "E83" = 1438
    "H83" = 627
    "B85" = 15319
    "C85" = 50
    "D85" = 14139
    "E85" = 980
    "G85" = 1
    "H85" = 200
    "B87" = 13987
    "C87" = 39
    "D87" = 9922
    "E87" = 3775
    "B95" = 9722
    "C95" = 73
    "D95" = 8886
    "E95" = 775
    "B96" = 9397
    "C96" = 6
    "D96" = 9115
    "E96" = 154
    "B103" = 8360
    "C103" = 24
    "E103" = 2213
    "G103" = 2
    "H103" = 214
    "B104" = 8291
    "C104" = 30
    "E104" = 605
    "B112" = 4879
    "C112" = 21
    "D112" = 4510
    "E112" = 275
    "E121" = 817
    "G121" = 1
    "H121" = 77
    "B128" = 3165
    "C128" = 43
    "D128" = 2483
    "E128" = 547
    "H128" = 135
    "B129" = 3150
    "D129" = 1315
    "E129" = 1736
    "H129" = 99
    "D130" = 2925
    "E130" = 184
    "B132" = 3083
    "C132" = 43
    "D132" = 1954
    "E132" = 1043
    "B147" = 2039
    "C147" = 25
    "D147" = 1627
    "E147" = 398
    "B183" = 315
    "C183" = 3
    "D183" = 273
    "E183" = 42
    "B218" = 6
    "C218" = 1
    "E218" = 1
}

foreach ($addr in $values.Keys) {
    $ws.Range($addr).Value = $values[$addr]
}
